$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New review row (row 5) - same appid/keyword as the row above it.
$ws.Range("A5").Value = "com.singleton.helix"
$ws.Range("B5").Value = "helix jump"
$ws.Range("C5").Value = "shmulmaor2@gmail.com"
$ws.Range("D5").Value = "vikicrestina@gmail.com"
$ws.Range("E5").Value = "27/5/2019 15:57"
$ws.Range("F5").Value = "Awesome game and awesome graphics. One of the best jump games and helix maze ever…"

# Hyperlink the two e-mail address cells (mirrors rows 2-4 above them),
# then restore the pre-existing cell formatting: Excel auto-applies its
# built-in "Hyperlink" style when a hyperlink is added to a cell, but the
# source cells here keep their original column formatting instead.
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:shmulmaor2@gmail.com", "", "", "shmulmaor2@gmail.com")
$ws.Range("C5").Font.Name = $ws.Range("C4").Font.Name
$ws.Range("C5").Font.Size = $ws.Range("C4").Font.Size
$ws.Range("C5").Font.Color = $ws.Range("C4").Font.Color
$ws.Range("C5").Font.Underline = $ws.Range("C4").Font.Underline

$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Range("D5").Font.Name = $ws.Range("D4").Font.Name
$ws.Range("D5").Font.Size = $ws.Range("D4").Font.Size
$ws.Range("D5").Font.Color = $ws.Range("D4").Font.Color
$ws.Range("D5").Font.Underline = $ws.Range("D4").Font.Underline

# Row 5 was previously a blank, taller spacer row (13.8pt); now that it
# holds real data it reverts to the sheet's normal row height (12.8pt).
$ws.Rows.Item(5).AutoFit()

# Selection moves to F5, the review cell of the newly added row.
$ws.Range("F5").Select()
